$wb = $excel.ActiveWorkbook

# --- Rename the two worksheets ---
$wsProcesses = $wb.Worksheets.Item("BaseIndustries")
$wsProcesses.Name = "BaseProcesses"

$wsFlows = $wb.Worksheets.Item("BaseProducts")
$wsFlows.Name = "BaseFlows"

# --- Add a "Type" column to the BaseFlows sheet, classifying each flow ---
$wsFlows.Range("C1").Value = "Type"
$wsFlows.Range("C1").Font.Bold = $true

$types = @(
    "Primary",    # 2  Iron ore
    "Primary",    # 3  Sponge iron
    "Primary",    # 4  Pig iron
    "Primary",    # 5  Liquid steel (OBF/OHF)
    "Primary",    # 6  Liquid steel (EAF)
    "Secondary",  # 7  Ingots
    "Secondary",  # 8  Slabs
    "Secondary",  # 9  Billets & blooms
    "Finished",   # 10 Flat rolled products
    "Finished",   # 11 Long rolled products
    "Recycling",  # 12 Forming & fabrication scrap
    "Recycling",  # 13 Scrap steel
    "Final",      # 14 Fabricated metal products, except machinery and equipment
    "Final",      # 15 Machinery and equipment n.e.c.
    "Final",      # 16 Office machinery and computers
    "Final",      # 17 Electrical machinery and apparatus n.e.c.
    "Final",      # 18 Radio, television and communication equipment and apparatus
    "Final",      # 19 Medical, precision and optical instruments, watches and clocks
    "Final",      # 20 Motor vehicles, trailers and semi-trailers
    "Final",      # 21 Other transport equipment
    "Final",      # 22 Furniture; other manufactured goods n.e.c.
    "Final"       # 23 Construction work
)

for ($i = 0; $i -lt $types.Length; $i++) {
    $row = $i + 2
    $wsFlows.Cells.Item($row, 3).Value = $types[$i]
}

# --- Activate BaseFlows and select C10, matching the saved UI state ---
$wsFlows.Activate()
$null = $wsFlows.Range("C10").Select()
